$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.3556784418185893
$ws.Range("D2").Value = 0.7254691660115968

$ws.Range("C3").Value = 0.3137912283741482
$ws.Range("D3").Value = 0.7566355952458692

$ws.Range("C4").Value = 2.536608715189334
$ws.Range("D4").Value = 0.01879816877881169

$ws.Range("C5").Value = 2.531968556851101
$ws.Range("D5").Value = 0.01899059516874302

$ws.Range("C6").Value = 0.5537539434153976
$ws.Range("D6").Value = 0.5853321147798374

$ws.Range("C7").Value = 3.158671669113759
$ws.Range("D7").Value = 0.00455395711800044

$ws.Range("C8").Value = 3.029454063033417
$ws.Range("D8").Value = 0.006158931445839011

$ws.Range("C9").Value = 2.089208695614153
$ws.Range("D9").Value = 0.04846939818333618

$ws.Range("C10").Value = 2.668123875305132
$ws.Range("D10").Value = 0.01404640878379193

$ws.Range("C11").Value = -0.4911625614382165
$ws.Range("D11").Value = 0.6281755534610189
